$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PINS")

# Row 12 (Accounts Payable) updates
$ws.Range("B12").Value = 49000000.0
$ws.Range("C12").Value = 49000000.0
$ws.Range("D12").Value = 42000000.0
$ws.Range("E12").Value = 44000000.0
$ws.Range("F12").Value = 38000000.0

# Row 37 (Net Debt) update
$ws.Range("G37").Value = -1493345000.0
